# Populate the rule-number column (A) for rows 3-5 of the "TripType"
# decision table. These cells were previously empty numeric cells
# (s="3" t="n"); they become text cells "1", "2", "3" while keeping the
# same visual style (border/font) as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TripType")

# Force the values to be stored as text (not auto-converted to numbers)
# by pre-formatting the target cells as Text, matching how the rule
# index labels are represented elsewhere in the DMN rule tables.
$rng = $ws.Range("A3:A5")
$rng.NumberFormat = "@"

$ws.Range("A3").Value = "1"
$ws.Range("A4").Value = "2"
$ws.Range("A5").Value = "3"
